$d = $word.ActiveDocument

# The document holds several "Test_Case N:" blocks, each a JSON object
# followed by a "Price: <value>" line. This edit swaps the JSON-body
# content (everything between the opening "{" and the trailing blank
# paragraph, including the Price line) of "Test_Case 1" and
# "Test_Case 2" while leaving the "Test_Case 1:" / "Test_Case 2:"
# header paragraphs themselves untouched.
#
# Paragraph numbers (1-based, Document.Paragraphs) in the original doc:
#   1  Test_Case 1:
#   2  {
#   3    "bedrooms": 3,
#   4    "bathrooms": 1.75,
#   5    "sqft_living": 1510,
#   6    "sqft_lot": 6360,
#   7    "floors": 1.0,
#   8    "waterfront": 0,
#   9    "view": 0,
#   10   "condition": 4,
#   11   "sqft_above": 1510,
#   12   "sqft_basement": 0,
#   13   "yr_built": 1954,
#   14   "yr_renovated": 1979,
#   15   "city": "Seattle"
#   16 }
#   17 Price: 308166.666667
#   18 (blank)
#   19 Test_Case 2:
#   20 {
#   21   "bedrooms": 3,
#   22   "bathrooms": 2.5,
#   23   "sqft_living": 1460,
#   24   "sqft_lot": 7573,
#   25   "floors": 2,
#   26   "waterfront": 0,
#   27   "view": 0,
#   28   "condition": 3,
#   29   "sqft_above": 1460,
#   30   "sqft_basement": 0,
#   31   "yr_built": 1983,
#   32   "yr_renovated": 2009,
#   33   "city": "Bellevue"
#   34 }
#   35 Price: 534333.333333
#
# Swap the values line-by-line between the two blocks (the surrounding
# "{", "}", "bedrooms" and blank-line paragraphs are identical between
# the two blocks, so only these lines actually change text).

$pairs = @(
    @(4,  '  "bathrooms": 1.75,',    22, '  "bathrooms": 2.5,'),
    @(5,  '  "sqft_living": 1510,',  23, '  "sqft_living": 1460,'),
    @(6,  '  "sqft_lot": 6360,',     24, '  "sqft_lot": 7573,'),
    @(7,  '  "floors": 1.0,',        25, '  "floors": 2,'),
    @(10, '  "condition": 4,',       28, '  "condition": 3,'),
    @(11, '  "sqft_above": 1510,',   29, '  "sqft_above": 1460,'),
    @(13, '  "yr_built": 1954,',     31, '  "yr_built": 1983,'),
    @(14, '  "yr_renovated": 1979,', 32, '  "yr_renovated": 2009,'),
    @(15, '  "city": "Seattle"',     33, '  "city": "Bellevue"')
)

foreach ($pair in $pairs) {
    $idx1 = $pair[0]
    $val1 = $pair[1]
    $idx2 = $pair[2]
    $val2 = $pair[3]

    $d.Paragraphs($idx1).Range.Text = $val2
    $d.Paragraphs($idx2).Range.Text = $val1
}

# "Price: " lines keep the label run and only swap the numeric run.
$price1 = $d.Paragraphs(17).Range
$price1.Find.Execute("308166.666667", $true, $false, $false, $false, $false, $true, 1, $false, "534333.333333", 2)

$price2 = $d.Paragraphs(35).Range
$price2.Find.Execute("534333.333333", $true, $false, $false, $false, $false, $true, 1, $false, "308166.666667", 2)
